# Adds a new "2023" data column (U) to the table, mirroring the formatting
# of column T (2022), and restores the selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, alignment) of the
# last existing data column (T, year 2022) onto the new column (U) before
# writing the 2023 values into it.
$ws.Range("T4:T14").Copy()
$ws.Range("U4:U14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 4: header year ---
$ws.Range("U4").Value = 2023

# --- Data rows 5-14 (2023 values) ---
$ws.Range("U5").Value = 3.3
$ws.Range("U6").Value = 1
$ws.Range("U7").Value = 1.6
$ws.Range("U8").Value = 9.1999999999999993
$ws.Range("U9").Value = 6.1
$ws.Range("U10").Value = 1.5
$ws.Range("U11").Value = 4
$ws.Range("U12").Value = 4.4000000000000004
$ws.Range("U13").Value = 4.7
$ws.Range("U14").Value = 0.5

# --- Selection moves back to B1 ---
$ws.Range("B1").Select()
